$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.541.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.29%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.942.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.36%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.41%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'57.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.89%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.40%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.227.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.22%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.99%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.812"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.42%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'13.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.89%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.26%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.940.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.55%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'36.457.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.26%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0865"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'228.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.30%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.70%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -6.47%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.29%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.44%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'160.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +9.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.57%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.74%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -5.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.59%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.58%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.04%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.28%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +8.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.83%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.57%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.49%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.12%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'15.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.340.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.80%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'86.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.64%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.12%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.34%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.119.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.15%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'43.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.47%  "
$ws.Range("E51").Style = "Normal"
